$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.957.32"
$ws.Range("E2").Value = "  +3.71%  "

# Row 3
$ws.Range("D3").Value = "3.790.16"
$ws.Range("E3").Value = "  +6.63%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "427.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.31%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.59%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.742"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.90%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.153"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.34%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000315"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.50%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.16%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +14.79%  "

# Row 14
$ws.Range("D14").Value = "4.379.75"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.37%  "

# Row 16
$ws.Range("E16").Value = "  +1.16%  "

# Row 17
$ws.Range("D17").Value = "3.802.07"
$ws.Range("E17").Value = "  +7.23%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.21%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +11.88%  "

# Row 20
$ws.Range("D20").Value = "66.129.59"
$ws.Range("E20").Value = "  +4.01%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "406.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.17%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.86%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.94%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.11%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "36.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.20%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +47.12%  "

# Row 27
$ws.Range("E27").Value = "  +9.89%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.01%  "

# Row 29
$ws.Range("E29").Value = "  -0.54%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +16.58%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "704.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.81%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.131"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +16.12%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.14%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.96%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +41.95%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.150"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.38%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "56.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.33%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0475"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.78%  "

# Row 40
$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +50.96%  "

# Row 41
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.142"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.87%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.76%  "

# Row 43
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D43").Value = "0.0₃0677"
$ws.Range("E43").Value = "  +3.17%  "

# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.43%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +9.52%  "

# Row 46
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.321"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +16.43%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.18%  "

# Row 48
$ws.Range("E48").Value = "  +7.86%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.57%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.57%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.55%  "
